# Update the loading_percent results for Case 3_22 (380 kV case)
# New simulation results replace the previous B2:N25 data block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,13
$arr[0,0] = 20.48249366587727
$arr[0,1] = 8.149653998291772
$arr[0,2] = 7.286513955439299
$arr[0,3] = 0
$arr[0,4] = 41.22157376090607
$arr[0,5] = 49.26234733102064
$arr[0,6] = 19.29036350191045
$arr[0,7] = 0
$arr[0,8] = 10.38853773159663
$arr[0,9] = 0
$arr[0,10] = 11.55922703570403
$arr[0,11] = 18.3671174808787
$arr[0,12] = 20.38727722585705
$arr[1,0] = 20.15935192383868
$arr[1,1] = 7.81519558711231
$arr[1,2] = 7.288029373004988
$arr[1,3] = 0
$arr[1,4] = 41.24672159956441
$arr[1,5] = 49.19256948707853
$arr[1,6] = 19.32915751996078
$arr[1,7] = 0
$arr[1,8] = 10.40406491697725
$arr[1,9] = 0
$arr[1,10] = 11.56899560773879
$arr[1,11] = 18.31263633060843
$arr[1,12] = 20.45854149768611
$arr[2,0] = 19.96299272784443
$arr[2,1] = 7.600756900924434
$arr[2,2] = 7.289212461993455
$arr[2,3] = 0
$arr[2,4] = 41.27210800648384
$arr[2,5] = 49.1643786576352
$arr[2,6] = 19.35672791430497
$arr[2,7] = 0
$arr[2,8] = 10.41408632593985
$arr[2,9] = 0
$arr[2,10] = 11.57631759686715
$arr[2,11] = 18.28222124452145
$arr[2,12] = 20.50424924492996
$arr[3,0] = 19.88359463691867
$arr[3,1] = 7.511159297351119
$arr[3,2] = 7.289758279655384
$arr[3,3] = 0
$arr[3,4] = 41.28494993176133
$arr[3,5] = 49.15657788198768
$arr[3,6] = 19.36890439373235
$arr[3,7] = 0
$arr[3,8] = 10.41829313076015
$arr[3,9] = 0
$arr[3,10] = 11.57963488269247
$arr[3,11] = 18.2705990431955
$arr[3,12] = 20.52336782579141
$arr[4,0] = 19.87045111781559
$arr[4,1] = 7.496150356253235
$arr[4,2] = 7.289852764744981
$arr[4,3] = 0
$arr[4,4] = 41.28723297747268
$arr[4,5] = 49.1555052716143
$arr[4,6] = 19.37098308585619
$arr[4,7] = 0
$arr[4,8] = 10.41899910737034
$arr[4,9] = 0
$arr[4,10] = 11.58020587495717
$arr[4,11] = 18.26871606214359
$arr[4,12] = 20.52657223417498
$arr[5,0] = 19.96191929034347
$arr[5,1] = 7.599557411586989
$arr[5,2] = 7.289219564917905
$arr[5,3] = 0
$arr[5,4] = 41.27227109433849
$arr[5,5] = 49.16425852367354
$arr[5,6] = 19.35688832177152
$arr[5,7] = 0
$arr[5,8] = 10.41414256186441
$arr[5,9] = 0
$arr[5,10] = 11.57636098382414
$arr[5,11] = 18.28206136602585
$arr[5,12] = 20.50450508931028
$arr[6,0] = 20.37071118428464
$arr[6,1] = 8.036256857754545
$arr[6,2] = 7.286984146067656
$arr[6,3] = 0
$arr[6,4] = 41.22817840849294
$arr[6,5] = 49.23524933213469
$arr[6,6] = 19.30296016398378
$arr[6,7] = 0
$arr[6,8] = 10.39379051862424
$arr[6,9] = 0
$arr[6,10] = 11.5623207884402
$arr[6,11] = 18.34770752444822
$arr[6,12] = 20.41144517964282
$arr[7,0] = 21.18405977975633
$arr[7,1] = 8.817942065546699
$arr[7,2] = 7.284597503045704
$arr[7,3] = 0
$arr[7,4] = 41.22076211583888
$arr[7,5] = 49.49045167831198
$arr[7,6] = 19.2270441564825
$arr[7,7] = 0
$arr[7,8] = 10.35773233895067
$arr[7,9] = 0
$arr[7,10] = 11.54526778616494
$arr[7,11] = 18.50012474098387
$arr[7,12] = 20.24435976467052
$arr[8,0] = 21.78272441874283
$arr[8,1] = 9.343792552525546
$arr[8,2] = 7.284052087835819
$arr[8,3] = 0
$arr[8,4] = 41.26361822676712
$arr[8,5] = 49.74805349864586
$arr[8,6] = 19.18955554729876
$arr[8,7] = 0
$arr[8,8] = 10.3335645129102
$arr[8,9] = 0
$arr[8,10] = 11.53909172704653
$arr[8,11] = 18.62595614615731
$arr[8,12] = 20.13088338330019
$arr[9,0] = 22.05407125614236
$arr[9,1] = 9.572024725388344
$arr[9,2] = 7.284064332047194
$arr[9,3] = 0
$arr[9,4] = 41.29360092545522
$arr[9,5] = 49.88025857026089
$arr[9,6] = 19.17648926349306
$arr[9,7] = 0
$arr[9,8] = 10.32306943815587
$arr[9,9] = 0
$arr[9,10] = 11.53765299503198
$arr[9,11] = 18.68606880858786
$arr[9,12] = 20.08125183417962
$arr[10,0] = 22.15658177593096
$arr[10,1] = 9.656840239129291
$arr[10,2] = 7.284106236760477
$arr[10,3] = 0
$arr[10,4] = 41.30645994809111
$arr[10,5] = 49.93245762546193
$arr[10,6] = 19.17211587991071
$arr[10,7] = 0
$arr[10,8] = 10.31916659234486
$arr[10,9] = 0
$arr[10,10] = 11.53730448168308
$arr[10,11] = 18.70923208234463
$arr[10,12] = 20.06274198204263
$arr[11,0] = 22.13451665118573
$arr[11,1] = 9.638645808794402
$arr[11,2] = 7.284095557042869
$arr[11,3] = 0
$arr[11,4] = 41.30362364216273
$arr[11,5] = 49.92112105399524
$arr[11,6] = 19.17303219771497
$arr[11,7] = 0
$arr[11,8] = 10.3200039694215
$arr[11,9] = 0
$arr[11,10] = 11.53737082314737
$arr[11,11] = 18.70422588320201
$arr[11,12] = 20.06671577935514
$arr[12,0] = 22.06251028286584
$arr[12,1] = 9.579035002083227
$arr[12,2] = 7.284067033641516
$arr[12,3] = 0
$arr[12,4] = 41.29462869923003
$arr[12,5] = 49.884510362958
$arr[12,6] = 19.176117942253
$arr[12,7] = 0
$arr[12,8] = 10.32274691945018
$arr[12,9] = 0
$arr[12,10] = 11.53762039326814
$arr[12,11] = 18.68796652652349
$arr[12,12] = 20.07972332371245
$arr[13,0] = 22.01836978835038
$arr[13,1] = 9.542310991091558
$arr[13,2] = 7.284054410578935
$arr[13,3] = 0
$arr[13,4] = 41.28931495379464
$arr[13,5] = 49.86236263302095
$arr[13,6] = 19.17808290164064
$arr[13,7] = 0
$arr[13,8] = 10.32443634525819
$arr[13,9] = 0
$arr[13,10] = 11.53779880159812
$arr[13,11] = 18.67805888685196
$arr[13,12] = 20.08772783000658
$arr[14,0] = 21.7649626953006
$arr[14,1] = 9.328653226911721
$arr[14,2] = 7.284056512073428
$arr[14,3] = 0
$arr[14,4] = 41.26186960744411
$arr[14,5] = 49.73971386484212
$arr[14,6] = 19.19048979992666
$arr[14,7] = 0
$arr[14,8] = 10.33426040231109
$arr[14,9] = 0
$arr[14,10] = 11.53921326688786
$arr[14,11] = 18.62208432398366
$arr[14,12] = 20.13416682234107
$arr[15,0] = 21.6091797581128
$arr[15,1] = 9.194744272231087
$arr[15,2] = 7.28412436310918
$arr[15,3] = 0
$arr[15,4] = 41.24771724539467
$arr[15,5] = 49.66830400006285
$arr[15,6] = 19.19912319295424
$arr[15,7] = 0
$arr[15,8] = 10.34041470166518
$arr[15,9] = 0
$arr[15,10] = 11.54043152416142
$arr[15,11] = 18.58847210304685
$arr[15,12] = 20.16316408186227
$arr[16,0] = 21.51949025744794
$arr[16,1] = 9.116691460168429
$arr[16,2] = 7.284187905438673
$arr[16,3] = 0
$arr[16,4] = 41.24056454771483
$arr[16,5] = 49.62864643217568
$arr[16,6] = 19.20446415625919
$arr[16,7] = 0
$arr[16,8] = 10.34400147956209
$arr[16,9] = 0
$arr[16,10] = 11.54126131057041
$arr[16,11] = 18.56941031905343
$arr[16,12] = 20.18002988403305
$arr[17,0] = 21.48911106650645
$arr[17,1] = 9.090087954952688
$arr[17,2] = 7.284213635659732
$arr[17,3] = 0
$arr[17,4] = 41.23831241153438
$arr[17,5] = 49.61546284865501
$arr[17,6] = 19.20633692861918
$arr[17,7] = 0
$arr[17,8] = 10.34522398160553
$arr[17,9] = 0
$arr[17,10] = 11.54156445507126
$arr[17,11] = 18.56300327041724
$arr[17,12] = 20.18577257869191
$arr[18,0] = 21.62577289016079
$arr[18,1] = 9.209106159348865
$arr[18,2] = 7.284114604165961
$arr[18,3] = 0
$arr[18,4] = 41.24912161568038
$arr[18,5] = 49.67575936627435
$arr[18,6] = 19.19816530616687
$arr[18,7] = 0
$arr[18,8] = 10.33975470543523
$arr[18,9] = 0
$arr[18,10] = 11.54028848534108
$arr[18,11] = 18.59202221817212
$arr[18,12] = 20.16005789567901
$arr[19,0] = 22.08366766535071
$arr[19,1] = 9.59658809462282
$arr[19,2] = 7.2840744015713
$arr[19,3] = 0
$arr[19,4] = 41.29722991013147
$arr[19,5] = 49.89520605363118
$arr[19,6] = 19.17519598348878
$arr[19,7] = 0
$arr[19,8] = 10.3219393131771
$arr[19,9] = 0
$arr[19,10] = 11.53754176760704
$arr[19,11] = 18.69273154637438
$arr[19,12] = 20.07589498255052
$arr[20,0] = 22.38147029766217
$arr[20,1] = 9.840428162187651
$arr[20,2] = 7.284265274809373
$arr[20,3] = 0
$arr[20,4] = 41.33744276103733
$arr[20,5] = 50.05106289043045
$arr[20,6] = 19.16353341556762
$arr[20,7] = 0
$arr[20,8] = 10.31071200643424
$arr[20,9] = 0
$arr[20,10] = 11.53689041995562
$arr[20,11] = 18.76087551138777
$arr[20,12] = 20.02254745372656
$arr[21,0] = 22.22269301553345
$arr[21,1] = 9.711155674710954
$arr[21,2] = 7.284143589449994
$arr[21,3] = 0
$arr[21,4] = 41.31517912918211
$arr[21,5] = 49.96675018303574
$arr[21,6] = 19.16945117007734
$arr[21,7] = 0
$arr[21,8] = 10.31666627246179
$arr[21,9] = 0
$arr[21,10] = 11.53713368190702
$arr[21,11] = 18.72429752312546
$arr[21,12] = 20.05086885095866
$arr[22,0] = 21.6182715255284
$arr[22,1] = 9.20261647055206
$arr[22,2] = 7.284118939774159
$arr[22,3] = 0
$arr[22,4] = 41.24848363544644
$arr[22,5] = 49.67238444199734
$arr[22,6] = 19.19859719087538
$arr[22,7] = 0
$arr[22,8] = 10.34005293835044
$arr[22,9] = 0
$arr[22,10] = 11.54035275008019
$arr[22,11] = 18.59041639298133
$arr[22,12] = 20.16146159520586
$arr[23,0] = 20.96340459204995
$arr[23,1] = 8.614802870325466
$arr[23,2] = 7.28503032984482
$arr[23,3] = 0
$arr[23,4] = 41.21428860399121
$arr[23,5] = 49.40904276386091
$arr[23,6] = 19.24437643515141
$arr[23,7] = 0
$arr[23,8] = 10.36707715132736
$arr[23,9] = 0
$arr[23,10] = 11.54876266659218
$arr[23,11] = 18.45641428862852
$arr[23,12] = 20.28792300439884

$ws.Range("B2:N25").Value = $arr
